$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Update the Usuario test value on Planilha1 (B2): usertest99 -> usertest206
$ws.Range("B2").Value = "usertest206"

# Move the saved cell selection on Planilha1 from F4 to F6
$ws.Activate()
$ws.Range("F6").Select()
